$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '62.800.14'
$ws.Range("E2").Value = '  +1.45%  '

# Row 3
$ws.Range("D3").Value = '2.442.45'
$ws.Range("E3").Value = '  +1.82%  '

# Row 4
$ws.Range("E4").Value = '  +0.16%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '567.85'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.23%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.04'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.63%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.05%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.534'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.05%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.111'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.83%  '

# Row 10
$ws.Range("E10").Value = '  +0.32%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.29'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.74%  '

# Row 12
$ws.Range("E12").Value = '  +2.14%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '27.01'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +5.89%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000182'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +6.15%  '

# Row 15
$ws.Range("D15").Value = '2.799.35'
$ws.Range("E15").Value = '  -1.16%  '

# Row 16
$ws.Range("D16").Value = '62.579.88'
$ws.Range("E16").Value = '  +1.32%  '

# Row 17
$ws.Range("D17").Value = '2.434.67'
$ws.Range("E17").Value = '  +1.63%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.27'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.66%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.94'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.27%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '324.40'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.28%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.17'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.13%  '

# Row 22
$ws.Range("E22").Value = '  -0.02%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.85'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +6.88%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '67.30'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.27%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '8.66'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.03%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '586.06'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +4.61%  '

# Row 27
$ws.Range("E27").Value = '  +9.51%  '

# Row 28
$ws.Range("D28").Value = '2.560.31'
$ws.Range("E28").Value = '  +1.59%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.46'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.90%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.999'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.48%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.45'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +4.65%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.147'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.61%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.87'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.58%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.54'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.95%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.87'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.59%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.999'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.10%  '

# Row 37
$ws.Range("E37").Value = '  +1.62%  '

# Row 38
$ws.Range("B38").Value = 'EthereumClassic'
$ws.Range("C38").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '18.83'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.70%  '

# Row 39
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.41'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.28%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '148.82'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.28%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.82'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.33%  '

# Row 42
$ws.Range("E42").Value = '  +0.35%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.46'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +10.17%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '148.97'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.29%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.68'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.51%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0537'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.57%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '20.62'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +4.69%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.602'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.04%  '

# Row 49
$ws.Range("E49").Value = '  +3.43%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0923'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.86%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.10'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +4.21%  '
